$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: remove Homework 0 / Brainstorm 3 topics (moved to row 4)
$ws.Range("F3").Value = ""
$ws.Range("H3").Value = ""

# Row 4: add Homework 0 / Brainstorm 3 topics
$ws.Range("F4").Value = "Homework 0"
$ws.Range("H4").Value = "Brainstorm 3 topics"

# Row 6: remove Topic and Timeline
$ws.Range("H6").Value = ""

# Row 7: remove Homework 1 Due
$ws.Range("F7").Value = ""

# Row 8: add Homework 1
$ws.Range("F8").Value = "Homework 1"

# Row 9: Outline -> Choose topic
$ws.Range("H9").Value = "Choose topic"

# Row 12: Homework 2 Due -> Homework 2; add Outline
$ws.Range("F12").Value = "Homework 2"
$ws.Range("H12").Value = "Outline"

# Row 13: remove Outline Feedback
$ws.Range("H13").Value = ""

# Row 14: add Outline Feedback
$ws.Range("H14").Value = "Outline Feedback"

# Row 16: add Homework 3
$ws.Range("F16").Value = "Homework 3"

# Row 18: remove Homework 3 Due
$ws.Range("F18").Value = ""

# Row 19: remove First draft
$ws.Range("H19").Value = ""

# Row 20: add Homework 4 and First draft
$ws.Range("F20").Value = "Homework 4"
$ws.Range("H20").Value = "First draft"

# Row 24: Homework 4 Due -> Homework 5
$ws.Range("F24").Value = "Homework 5"

# Row 27: remove article and discussion
$ws.Range("G27").Value = ""

# Selection moved to A20 as part of the class-presentation review pass
$ws.Range("A20").Select()
